function Set-CellText {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-CellText $ws 'D2' '65.674.80'
Set-CellText $ws 'E2' '  -0.11%  '

# Row 3
Set-CellText $ws 'D3' '3.298.95'
Set-CellText $ws 'E3' '  -0.34%  '

# Row 5
Set-CellText $ws 'D5' '186.73'
Set-CellText $ws 'E5' '  +3.71%  '

# Row 6
Set-CellText $ws 'D6' '555.26'
Set-CellText $ws 'E6' '  -0.15%  '

# Row 7
Set-CellText $ws 'E7' '  +0.03%  '

# Row 8
Set-CellText $ws 'D8' '0.582'
Set-CellText $ws 'E8' '  -0.53%  '

# Row 9
Set-CellText $ws 'D9' '3.290.38'
Set-CellText $ws 'E9' '  -0.39%  '

# Row 10
Set-CellText $ws 'D10' '0.182'
Set-CellText $ws 'E10' '  -0.31%  '

# Row 11
Set-CellText $ws 'D11' '0.583'
Set-CellText $ws 'E11' '  +0.65%  '

# Row 12
Set-CellText $ws 'D12' '47.17'
Set-CellText $ws 'E12' '  +0.57%  '

# Row 13
Set-CellText $ws 'D13' '0.0000268'
Set-CellText $ws 'E13' '  +2.94%  '

# Row 14
Set-CellText $ws 'D14' '8.68'
Set-CellText $ws 'E14' '  +2.59%  '

# Row 15
Set-CellText $ws 'D15' '3.832.62'
Set-CellText $ws 'E15' '  -0.40%  '

# Row 16
Set-CellText $ws 'D16' '604.46'
Set-CellText $ws 'E16' '  +1.71%  '

# Row 17
Set-CellText $ws 'D17' '65.852.18'
Set-CellText $ws 'E17' '  +0.17%  '

# Row 18
Set-CellText $ws 'B18' 'TRON'
Set-CellText $ws 'C18' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText $ws 'D18' '0.118'
Set-CellText $ws 'E18' '  +1.11%  '

# Row 19
Set-CellText $ws 'B19' 'Chainlink'
Set-CellText $ws 'C19' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText $ws 'D19' '17.94'
Set-CellText $ws 'E19' '  +0.41%  '

# Row 20
Set-CellText $ws 'D20' '3.312.08'
Set-CellText $ws 'E20' '  -0.10%  '

# Row 21
Set-CellText $ws 'D21' '11.05'
Set-CellText $ws 'E21' '  -1.97%  '

# Row 22
Set-CellText $ws 'D22' '0.905'
Set-CellText $ws 'E22' '  +0.96%  '

# Row 23
Set-CellText $ws 'D23' '18.56'
Set-CellText $ws 'E23' '  +11.45%  '

# Row 24
Set-CellText $ws 'D24' '5.07'
Set-CellText $ws 'E24' '  +0.86%  '

# Row 25
Set-CellText $ws 'D25' '99.92'
Set-CellText $ws 'E25' '  +0.89%  '

# Row 26
Set-CellText $ws 'D26' '3.95'
Set-CellText $ws 'E26' '  -0.37%  '

# Row 27
Set-CellText $ws 'D27' '2.76'
Set-CellText $ws 'E27' '  +5.20%  '

# Row 28
Set-CellText $ws 'D28' '5.92'
Set-CellText $ws 'E28' '  -1.14%  '

# Row 29
Set-CellText $ws 'D29' '9.54'
Set-CellText $ws 'E29' '  +4.20%  '

# Row 30
Set-CellText $ws 'D30' '8.68'
Set-CellText $ws 'E30' '  +1.22%  '

# Row 31
Set-CellText $ws 'D31' '30.23'
Set-CellText $ws 'E31' '  -0.53%  '

# Row 32
Set-CellText $ws 'D32' '6.72'
Set-CellText $ws 'E32' '  +9.04%  '

# Row 33
Set-CellText $ws 'D33' '580.43'
Set-CellText $ws 'E33' '  +10.74%  '

# Row 34
Set-CellText $ws 'E34' '  -0.18%  '

# Row 35
Set-CellText $ws 'D35' '11.10'
Set-CellText $ws 'E35' '  +1.80%  '

# Row 36
Set-CellText $ws 'E36' '  +1.49%  '

# Row 37
Set-CellText $ws 'B37' 'Dai'
Set-CellText $ws 'C37' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText $ws 'D37' '0.999'
Set-CellText $ws 'E37' '  +0.01%  '

# Row 38
Set-CellText $ws 'B38' 'OKB'
Set-CellText $ws 'C38' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText $ws 'D38' '57.02'
Set-CellText $ws 'E38' '  -0.74%  '

# Row 39
Set-CellText $ws 'B39' 'Maker'
Set-CellText $ws 'C39' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws 'D39' '3.698.22'
Set-CellText $ws 'E39' '  -0.73%  '

# Row 40
Set-CellText $ws 'B40' 'InjectiveProtocol'
Set-CellText $ws 'C40' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-CellText $ws 'D40' '33.88'
Set-CellText $ws 'E40' '  +7.66%  '

# Row 41
Set-CellText $ws 'D41' '0.130'
Set-CellText $ws 'E41' '  +5.65%  '

# Row 42
Set-CellText $ws 'B42' 'CoreDAO'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-CellText $ws 'D42' '3.45'
Set-CellText $ws 'E42' '  +16.00%  '

# Row 43
Set-CellText $ws 'B43' 'PEPE'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText $ws 'D43' '0.0₃0719'
Set-CellText $ws 'E43' '  +2.54%  '

# Row 44
Set-CellText $ws 'D44' '3.26'
Set-CellText $ws 'E44' '  -4.70%  '

# Row 45
Set-CellText $ws 'D45' '2.66'
Set-CellText $ws 'E45' '  +1.47%  '

# Row 46
Set-CellText $ws 'D46' '0.340'
Set-CellText $ws 'E46' '  +1.26%  '

# Row 47
Set-CellText $ws 'D47' '3.36'
Set-CellText $ws 'E47' '  +3.01%  '

# Row 48
Set-CellText $ws 'D48' '0.0420'
Set-CellText $ws 'E48' '  +2.72%  '

# Row 49
Set-CellText $ws 'D49' '0.129'
Set-CellText $ws 'E49' '  +0.97%  '

# Row 50
Set-CellText $ws 'D50' '2.58'
Set-CellText $ws 'E50' '  +0.03%  '

# Row 51
Set-CellText $ws 'D51' '1.00'
Set-CellText $ws 'E51' '  +0.22%  '
